$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 797.875
$ws.Range("J17").Value = 817.5714
$ws.Range("L17").Value = 2452.7142
$ws.Range("N17").Value = -2788.7142

$ws.Range("H20").Value = 1268.7142
$ws.Range("I20").Value = 1230.1666
$ws.Range("J20").Value = 1500
$ws.Range("K20").Value = 1230.1666
$ws.Range("L20").Value = 1500
$ws.Range("M20").Value = -1000.1666
$ws.Range("N20").Value = -1960

$ws.Range("H35").Value = 1268.7142
$ws.Range("I35").Value = 1230.1666
$ws.Range("J35").Value = 1500
$ws.Range("K35").Value = 1230.1666
$ws.Range("L35").Value = 1500
$ws.Range("M35").Value = -851.1666
$ws.Range("N35").Value = -2258

$ws.Range("H64").Value = 4513
$ws.Range("J64").Value = 4798.8
$ws.Range("L64").Value = 4798.8
$ws.Range("N64").Value = -5294.8

$ws.Range("H67").Value = 4513
$ws.Range("J67").Value = 4798.8
$ws.Range("L67").Value = 4798.8
$ws.Range("N67").Value = -6514.8

$ws.Range("H94").Value = 1692.75
$ws.Range("I94").Value = 1593.6666
$ws.Range("K94").Value = 1593.6666
$ws.Range("M94").Value = -1142.6666

$ws.Range("H132").Value = 3588.875
$ws.Range("I132").Value = 3673
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 11019
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -8489
$ws.Range("N132").Value = -14060

$ws.Range("H137").Value = 4349.1
$ws.Range("I137").Value = 2785.2856
$ws.Range("J137").Value = 7998
$ws.Range("K137").Value = 8355.856800000001
$ws.Range("L137").Value = 23994
$ws.Range("M137").Value = -5805.856800000001
$ws.Range("N137").Value = -29094

$ws.Range("H138").Value = 2383.3044
$ws.Range("I138").Value = 1744.5834
$ws.Range("J138").Value = 3080.0908
$ws.Range("K138").Value = 5233.7502
$ws.Range("L138").Value = 9240.2724
$ws.Range("M138").Value = -93.7502000000004
$ws.Range("N138").Value = -19520.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 962
$ws.Range("I2").Value = 1224.5
$ws.Range("J2").Value = 699.5
$ws.Range("K2").Value = 1224.5
$ws.Range("L2").Value = 699.5
$ws.Range("M2").Value = -1111.5
$ws.Range("N2").Value = -925.5

$ws.Range("H45").Value = 6120
$ws.Range("I45").Value = 8966.666999999999
$ws.Range("K45").Value = 8966.666999999999
$ws.Range("M45").Value = -8589.666999999999

$ws.Range("H116").Value = 962
$ws.Range("I116").Value = 1224.5
$ws.Range("J116").Value = 699.5
$ws.Range("K116").Value = 1224.5
$ws.Range("L116").Value = 699.5
$ws.Range("M116").Value = 1069.5
$ws.Range("N116").Value = -5287.5

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = $null

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 962
$ws.Range("I3").Value = 1224.5
$ws.Range("J3").Value = 699.5
$ws.Range("K3").Value = 1224.5
$ws.Range("L3").Value = 699.5
$ws.Range("M3").Value = -1110.5
$ws.Range("N3").Value = -927.5

$ws.Range("H88").Value = 21577.8
$ws.Range("J88").Value = 24972.5
$ws.Range("L88").Value = 24972.5
$ws.Range("N88").Value = -25784.5

$ws.Range("H91").Value = 21577.8
$ws.Range("J91").Value = 24972.5
$ws.Range("L91").Value = 24972.5
$ws.Range("N91").Value = -27780.5

$ws.Range("H94").Value = 649.3333
$ws.Range("I94").Value = 724.5
$ws.Range("K94").Value = 724.5
$ws.Range("M94").Value = -273.5

$ws.Range("H99").Value = 1991
$ws.Range("I99").Value = 1487.375
$ws.Range("J99").Value = 4005.5
$ws.Range("K99").Value = 1487.375
$ws.Range("L99").Value = 4005.5
$ws.Range("M99").Value = 10.625
$ws.Range("N99").Value = -7001.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 199142
$ws.Range("J9").Value = 199142
$ws.Range("L9").Value = 199142
$ws.Range("N9").Value = -199478

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = $null

$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").Value = $null

$ws.Range("H107").Value = 1241
$ws.Range("I107").Value = 1196.75
$ws.Range("K107").Value = 1196.75
$ws.Range("M107").Value = 723.25

$ws.Range("H132").Value = 6979.75
$ws.Range("I132").Value = 7489
$ws.Range("J132").Value = 6470.5
$ws.Range("K132").Value = 22467
$ws.Range("L132").Value = 19411.5
$ws.Range("M132").Value = -19937
$ws.Range("N132").Value = -24471.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1920
$ws.Range("J22").Value = 1920
$ws.Range("L22").Value = 5760
$ws.Range("N22").Value = -6098

$ws.Range("H27").Value = 1920
$ws.Range("J27").Value = 1920
$ws.Range("L27").Value = 5760
$ws.Range("N27").Value = -5964

$ws.Range("H70").Value = 13916.667
$ws.Range("I70").Value = 1750
$ws.Range("J70").Value = 20000
$ws.Range("K70").Value = 5250
$ws.Range("L70").Value = 60000
$ws.Range("M70").Value = -4935
$ws.Range("N70").Value = -60630

$ws.Range("H73").Value = 13916.667
$ws.Range("I73").Value = 1750
$ws.Range("J73").Value = 20000
$ws.Range("K73").Value = 5250
$ws.Range("L73").Value = 60000
$ws.Range("M73").Value = -4158
$ws.Range("N73").Value = -62184

$ws.Range("H116").Value = 2000
$ws.Range("J116").Value = 1000
$ws.Range("L116").Value = 3000
$ws.Range("N116").Value = -9884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 120.125
$ws.Range("I2").Value = 80.25
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 80.25
$ws.Range("L2").Value = 160
$ws.Range("M2").Value = 32.75
$ws.Range("N2").Value = -386

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("N113").Value = $null

$ws.Range("H122").Value = 2779.8333
$ws.Range("I122").Value = 2779.8333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8339.499899999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5889.499899999999
$ws.Range("N122").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3333.3333
$ws.Range("I22").Value = 3333.3333
$ws.Range("K22").Value = 3333.3333
$ws.Range("M22").Value = -3038.3333

$ws.Range("H27").Value = 3333.3333
$ws.Range("I27").Value = 3333.3333
$ws.Range("K27").Value = 3333.3333
$ws.Range("M27").Value = -3226.3333

$ws.Range("H55").Value = 747.2308
$ws.Range("J55").Value = 870.2857
$ws.Range("L55").Value = 870.2857
$ws.Range("N55").Value = -1216.2857

$ws.Range("H93").Value = 3898
$ws.Range("I93").Value = 3898
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 3898
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -2650
$ws.Range("N93").Value = $null

$ws.Range("H100").Value = 5619.6
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 30000
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = $null

$ws.Range("H132").Value = 2536.0908
$ws.Range("I132").Value = 2722
$ws.Range("J132").Value = 1699.5
$ws.Range("K132").Value = 8166
$ws.Range("L132").Value = 5098.5
$ws.Range("M132").Value = -5636
$ws.Range("N132").Value = -10158.5
